$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032618576973103
$ws.Range("D2").Value = 1.035584673880938
$ws.Range("E2").Value = 1.036259756061105
$ws.Range("F2").Value = 1.042237979673527
$ws.Range("I2").Value = 1.036790426248109
$ws.Range("J2").Value = 1.037747804950852
$ws.Range("K2").Value = 1.038380858729643
$ws.Range("L2").Value = 1.039054007935447
$ws.Range("M2").Value = 1.045015232607604
$ws.Range("N2").Value = 1.039221525865245
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033531309452188
$ws.Range("D3").Value = 1.036254382546084
$ws.Range("E3").Value = 1.037120419450005
$ws.Range("F3").Value = 1.043394940543585
$ws.Range("I3").Value = 1.037018912295743
$ws.Range("J3").Value = 1.038303094551655
$ws.Range("K3").Value = 1.038860363994916
$ws.Range("L3").Value = 1.03972410036771
$ws.Range("M3").Value = 1.045982077105118
$ws.Range("N3").Value = 1.039777604040948
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03412229751514
$ws.Range("D4").Value = 1.036688034300578
$ws.Range("E4").Value = 1.037678057693334
$ws.Range("F4").Value = 1.044144330654435
$ws.Range("I4").Value = 1.037165775952123
$ws.Range("J4").Value = 1.038662168026167
$ws.Range("K4").Value = 1.03917026027459
$ws.Range("L4").Value = 1.040157786149208
$ws.Range("M4").Value = 1.046607872162349
$ws.Range("N4").Value = 1.040137187440981
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034370841174559
$ws.Range("D5").Value = 1.036870413396332
$ws.Range("E5").Value = 1.037912662732144
$ws.Range("F5").Value = 1.044459555787061
$ws.Range("I5").Value = 1.037227281878971
$ws.Range("J5").Value = 1.038813065221411
$ws.Range("K5").Value = 1.039300449924562
$ws.Range("L5").Value = 1.040340128528706
$ws.Range("M5").Value = 1.046870999690087
$ws.Range("N5").Value = 1.040288298927549
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034412578134417
$ws.Range("D6").Value = 1.036901039829353
$ws.Range("E6").Value = 1.03795206411455
$ws.Range("F6").Value = 1.044512494113498
$ws.Range("I6").Value = 1.037237595169001
$ws.Range("J6").Value = 1.038838398141889
$ws.Range("K6").Value = 1.039322303994885
$ws.Range("L6").Value = 1.040370745813237
$ws.Range("M6").Value = 1.046915182489836
$ws.Range("N6").Value = 1.04031366782368
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03412561820618
$ws.Range("D7").Value = 1.036690470978016
$ws.Range("E7").Value = 1.037681191815448
$ws.Range("F7").Value = 1.044148541995019
$ws.Range("I7").Value = 1.037166598723437
$ws.Range("J7").Value = 1.038664184547982
$ws.Range("K7").Value = 1.03917200023264
$ws.Range("L7").Value = 1.040160222536161
$ws.Range("M7").Value = 1.046611387915518
$ws.Range("N7").Value = 1.040139206826488
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032926958089893
$ws.Range("D8").Value = 1.035810941071202
$ws.Range("E8").Value = 1.036550469105514
$ws.Range("F8").Value = 1.042628822729855
$ws.Range("I8").Value = 1.03686784726237
$ws.Range("J8").Value = 1.037935516068383
$ws.Range("K8").Value = 1.038542987083995
$ws.Range("L8").Value = 1.039280449361676
$ws.Range("M8").Value = 1.045341944368662
$ws.Range("N8").Value = 1.039409503554068
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030817784683126
$ws.Range("D9").Value = 1.034263490777815
$ws.Range("E9").Value = 1.034563646256555
$ws.Range("F9").Value = 1.039956705008889
$ws.Range("I9").Value = 1.036333906727899
$ws.Range("J9").Value = 1.036649736138475
$ws.Range("K9").Value = 1.037431743697016
$ws.Range("L9").Value = 1.037730916735177
$ws.Range("M9").Value = 1.043106429300585
$ws.Range("N9").Value = 1.038121897669178
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029413749592092
$ws.Range("D10").Value = 1.033233539234756
$ws.Range("E10").Value = 1.033242970036935
$ws.Range("F10").Value = 1.038179216001277
$ws.Range("I10").Value = 1.03597292793918
$ws.Range("J10").Value = 1.035791403471579
$ws.Range("K10").Value = 1.036689051254029
$ws.Range("L10").Value = 1.036698448161976
$ws.Range("M10").Value = 1.041617046519934
$ws.Range("N10").Value = 1.037262346071443
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028806290623346
$ws.Range("D11").Value = 1.032787972889885
$ws.Range("E11").Value = 1.032672035631992
$ws.Range("F11").Value = 1.037410475024403
$ws.Range("I11").Value = 1.035815435085012
$ws.Range("J11").Value = 1.035419473549106
$ws.Range("K11").Value = 1.036367024951294
$ws.Range("L11").Value = 1.036251520509911
$ws.Range("M11").Value = 1.040972357577841
$ws.Range("N11").Value = 1.036889887965819
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028580728826633
$ws.Range("D12").Value = 1.032622532423284
$ws.Range("E12").Value = 1.032460105598535
$ws.Range("F12").Value = 1.037125069114671
$ws.Range("I12").Value = 1.03575675733455
$ws.Range("J12").Value = 1.035281282940768
$ws.Range("K12").Value = 1.036247345354344
$ws.Range("L12").Value = 1.036085533305653
$ws.Range("M12").Value = 1.040732925098357
$ws.Range("N12").Value = 1.036751501110967
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028629109168276
$ws.Range("D13").Value = 1.032658017107925
$ws.Range("E13").Value = 1.032505558937315
$ws.Range("F13").Value = 1.037186283374254
$ws.Range("I13").Value = 1.035769351954921
$ws.Range("J13").Value = 1.035310927061045
$ws.Range("K13").Value = 1.036273019949454
$ws.Range("L13").Value = 1.03612113713163
$ws.Range("M13").Value = 1.040784282651377
$ws.Range("N13").Value = 1.036781187329295
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028787644069992
$ws.Range("D14").Value = 1.032774296243083
$ws.Range("E14").Value = 1.032654514558099
$ws.Range("F14").Value = 1.037386880453211
$ws.Range("I14").Value = 1.035810588390536
$ws.Range("J14").Value = 1.035408051467335
$ws.Range("K14").Value = 1.036357133513789
$ws.Range("L14").Value = 1.036237799509561
$ws.Range("M14").Value = 1.040952565326354
$ws.Range("N14").Value = 1.036878449663383
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028885332684324
$ws.Range("D15").Value = 1.032845948001698
$ws.Range("E15").Value = 1.032746309673976
$ws.Range("F15").Value = 1.0375104933148
$ws.Range("I15").Value = 1.035835971958766
$ws.Range("J15").Value = 1.03546788783683
$ws.Range("K15").Value = 1.036408950107289
$ws.Range("L15").Value = 1.036309681942965
$ws.Range("M15").Value = 1.04105625430492
$ws.Range("N15").Value = 1.036938371007386
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029454075190333
$ws.Range("D16").Value = 1.033263118734683
$ws.Range("E16").Value = 1.033280880703319
$ws.Range("F16").Value = 1.038230254262579
$ws.Range("I16").Value = 1.03598335525294
$ws.Range("J16").Value = 1.035816081688277
$ws.Range("K16").Value = 1.036710413966147
$ws.Range("L16").Value = 1.036728112262397
$ws.Range("M16").Value = 1.041659837121849
$ws.Range("N16").Value = 1.037287059334039
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029810966072438
$ws.Range("D17").Value = 1.033524909384454
$ws.Range("E17").Value = 1.033616452003941
$ws.Range("F17").Value = 1.038681988513382
$ws.Range("I17").Value = 1.036075487472306
$ws.Range("J17").Value = 1.036034423727532
$ws.Range("K17").Value = 1.036899398194115
$ws.Range("L17").Value = 1.036990620131696
$ws.Range("M17").Value = 1.042038508617731
$ws.Range("N17").Value = 1.037505711444034
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030019182434726
$ws.Range("D18").Value = 1.033677646869212
$ws.Range("E18").Value = 1.033812274442879
$ws.Range("F18").Value = 1.038945566564258
$ws.Range("I18").Value = 1.036129112115704
$ws.Range("J18").Value = 1.036161753210828
$ws.Range("K18").Value = 1.037009587420021
$ws.Range("L18").Value = 1.037143749791373
$ws.Range("M18").Value = 1.042259402864667
$ws.Range("N18").Value = 1.037633221749802
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030090186943096
$ws.Range("D19").Value = 1.033729733055713
$ws.Range("E19").Value = 1.033879059959304
$ws.Range("F19").Value = 1.039035454982844
$ws.Range("I19").Value = 1.036147377296677
$ws.Range("J19").Value = 1.036205164887281
$ws.Range("K19").Value = 1.037047151938724
$ws.Range("L19").Value = 1.037195965270473
$ws.Range("M19").Value = 1.042334725750139
$ws.Range("N19").Value = 1.037676695075814
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029772670093194
$ws.Range("D20").Value = 1.033496817638435
$ws.Range("E20").Value = 1.033580439126305
$ws.Range("F20").Value = 1.038633512486973
$ws.Range("I20").Value = 1.036065614399758
$ws.Range("J20").Value = 1.036011000349659
$ws.Range("K20").Value = 1.036879126321755
$ws.Range("L20").Value = 1.036962454143627
$ws.Range("M20").Value = 1.041997878518035
$ws.Range("N20").Value = 1.037482254802278
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028740957402135
$ws.Range("D21").Value = 1.032740053195197
$ws.Range("E21").Value = 1.032610646954836
$ws.Range("F21").Value = 1.037327805784115
$ws.Range("I21").Value = 1.035798450195589
$ws.Range("J21").Value = 1.035379451828833
$ws.Range("K21").Value = 1.036332365945298
$ws.Range("L21").Value = 1.036203444742926
$ws.Range("M21").Value = 1.04090300934183
$ws.Range("N21").Value = 1.036849809410113
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028092715954183
$ws.Range("D22").Value = 1.032264608780572
$ws.Range("E22").Value = 1.032001713147206
$ws.Range("F22").Value = 1.036507658723997
$ws.Range("I22").Value = 1.035629444422723
$ws.Range("J22").Value = 1.034982145332734
$ws.Range("K22").Value = 1.035988222091909
$ws.Range("L22").Value = 1.035726350962898
$ws.Range("M22").Value = 1.04021481665846
$ws.Range("N22").Value = 1.036451938693218
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028436319304879
$ws.Range("D23").Value = 1.032516616013349
$ws.Range("E23").Value = 1.032324443035098
$ws.Range("F23").Value = 1.036942358131339
$ws.Range("I23").Value = 1.035719134990431
$ws.Range("J23").Value = 1.035192786287736
$ws.Range("K23").Value = 1.03617069445109
$ws.Range("L23").Value = 1.035979255205866
$ws.Range("M23").Value = 1.040579622099834
$ws.Range("N23").Value = 1.036662878782536
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029789974249659
$ws.Range("D24").Value = 1.033509510967561
$ws.Range("E24").Value = 1.033596711519493
$ws.Range("F24").Value = 1.03865541644122
$ws.Range("I24").Value = 1.036070075970135
$ws.Range("J24").Value = 1.036021584445434
$ws.Range("K24").Value = 1.036888286438352
$ws.Range("L24").Value = 1.036975181100339
$ws.Range("M24").Value = 1.04201623744443
$ws.Range("N24").Value = 1.037492853928683
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031362693593843
$ws.Range("D25").Value = 1.034663252615268
$ws.Range("E25").Value = 1.035076610202938
$ws.Range("F25").Value = 1.040646820719614
$ws.Range("I25").Value = 1.036472829575779
$ws.Range("J25").Value = 1.036982346440947
$ws.Range("K25").Value = 1.037719358663395
$ws.Range("L25").Value = 1.038131414631824
$ws.Range("M25").Value = 1.043684195326604
$ws.Range("N25").Value = 1.038454980316431
